$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 0.1140306016422357
$ws.Range("J2").Value = 0.1140306016422357
$ws.Range("M2").Value = 179.7005413333333
$ws.Range("N2").Value = 539.101624
$ws.Range("O2").Value = 0.7012656334041908
$ws.Range("P2").Value = 0.7012656334041907
$ws.Range("Q2").Value = 42.07155063714045
$ws.Range("R2").Value = 378.643955734264
$ws.Range("S2").Value = 0.07996574208810341
$ws.Range("T2").Value = 0.07996574208810339
$ws.Range("I3").Value = 0.1140306016422357
$ws.Range("J3").Value = 0.1140306016422357
$ws.Range("O3").Value = 0.05908927597267952
$ws.Range("P3").Value = 0.05908927597267952
$ws.Range("S3").Value = 0.00673798568976875
$ws.Range("T3").Value = 0.006737985689768751
$ws.Range("I4").Value = 0.1140306016422357
$ws.Range("J4").Value = 0.1140306016422357
$ws.Range("M4").Value = 36.14947766666667
$ws.Range("N4").Value = 108.448433
$ws.Range("O4").Value = 0.1410701724382803
$ws.Range("P4").Value = 0.1410701724382803
$ws.Range("Q4").Value = 8.463327761145891
$ws.Range("R4").Value = 76.169949850313
$ws.Range("S4").Value = 0.01608631663691104
$ws.Range("T4").Value = 0.01608631663691104
$ws.Range("I5").Value = 0.1140306016422357
$ws.Range("J5").Value = 0.1140306016422357
$ws.Range("M5").Value = 25.25999466666667
$ws.Range("N5").Value = 75.779984
$ws.Range("O5").Value = 0.09857491818484938
$ws.Range("P5").Value = 0.09857491818484938
$ws.Range("Q5").Value = 5.913878371358223
$ws.Range("R5").Value = 53.224905342224
$ws.Range("S5").Value = 0.01124055722745254
$ws.Range("T5").Value = 0.01124055722745254
$ws.Range("G6").Value = 1.819015666666667
$ws.Range("H6").Value = 5.457047
$ws.Range("I6").Value = 0.8859693983577642
$ws.Range("J6").Value = 0.8859693983577643
$ws.Range("M6").Value = 179.7005413333333
$ws.Range("N6").Value = 539.101624
$ws.Range("O6").Value = 0.7012656334041908
$ws.Range("P6").Value = 0.7012656334041907
$ws.Range("Q6").Value = 326.8780999938143
$ws.Range("R6").Value = 2941.902899944328
$ws.Range("S6").Value = 0.6212998913160873
$ws.Range("T6").Value = 0.6212998913160873
$ws.Range("G7").Value = 1.819015666666667
$ws.Range("H7").Value = 5.457047
$ws.Range("I7").Value = 0.8859693983577642
$ws.Range("J7").Value = 0.8859693983577643
$ws.Range("O7").Value = 0.05908927597267952
$ws.Range("P7").Value = 0.05908927597267952
$ws.Range("Q7").Value = 27.54304409043667
$ws.Range("R7").Value = 247.88739681393
$ws.Range("S7").Value = 0.05235129028291077
$ws.Range("T7").Value = 0.05235129028291078
$ws.Range("G8").Value = 1.819015666666667
$ws.Range("H8").Value = 5.457047
$ws.Range("I8").Value = 0.8859693983577642
$ws.Range("J8").Value = 0.8859693983577643
$ws.Range("M8").Value = 36.14947766666667
$ws.Range("N8").Value = 108.448433
$ws.Range("O8").Value = 0.1410701724382803
$ws.Range("P8").Value = 0.1410701724382803
$ws.Range("Q8").Value = 65.75646621748345
$ws.Range("R8").Value = 591.808195957351
$ws.Range("S8").Value = 0.1249838558013692
$ws.Range("T8").Value = 0.1249838558013693
$ws.Range("G9").Value = 1.819015666666667
$ws.Range("H9").Value = 5.457047
$ws.Range("I9").Value = 0.8859693983577642
$ws.Range("J9").Value = 0.8859693983577643
$ws.Range("M9").Value = 25.25999466666667
$ws.Range("N9").Value = 75.779984
$ws.Range("O9").Value = 0.09857491818484938
$ws.Range("P9").Value = 0.09857491818484938
$ws.Range("Q9").Value = 45.94832603858312
$ws.Range("R9").Value = 413.534934347248
$ws.Range("S9").Value = 0.08733436095739684
$ws.Range("T9").Value = 0.08733436095739686
